# Update command-line documentation table: add "dfa" and "nfa" topics to the
# --show <topic> help text, and change the long-form switch to use "=" syntax.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Long form switch text changes from "--show <topic>" to "--show=<topic>"
$ws.Range("B12").Value = "--show=<topic>"

# Description text gains two new bullet points (dfa, nfa) and switches the
# bullet separator from a tab to a single space.
$newDescription = "Show information on a range of topics:" + [char]10 + `
    [char]8226 + " dfa " + [char]8211 + " the Deterministic Finite Automata for the expression parser" + [char]10 + `
    [char]8226 + " distribution " + [char]8211 + " the licence rules on distributing the software" + [char]10 + `
    [char]8226 + " grammar " + [char]8211 + " the full list of grammar rules being applied" + [char]10 + `
    [char]8226 + " instructions " + [char]8211 + " the full list of instructions available" + [char]10 + `
    [char]8226 + " nfa " + [char]8211 + " the Non-deterministic Finite Automata for the expression parser" + [char]10 + `
    [char]8226 + " operators " + [char]8211 + " the operators which are used and their priority" + [char]10 + `
    [char]8226 + " reserved " + [char]8211 + " the reserved words used for the chosen processor and grammar" + [char]10 + `
    [char]8226 + " version " + [char]8211 + " the version of the software and other key details" + [char]10 + `
    [char]8226 + " warranty " + [char]8211 + " the warranty available with the software"

$ws.Range("D12").Value = $newDescription

# Reset the selection back to the top-left cell so the saved sheet view no
# longer records a stale selection at C13.
$ws.Range("A1").Select()

$wb.Save()
